$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.885.84'
$ws.Range("E2").Value = '  +1.55%  '

$ws.Range("D3").Value = '3.136.98'
$ws.Range("E3").Value = '  +3.28%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.94%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.128.45'
$ws.Range("E8").Value = '  +3.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.94%  '

$ws.Range("E11").Value = '  +2.43%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.57%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000222'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '

$ws.Range("D15").Value = '3.641.00'
$ws.Range("E15").Value = '  +3.34%  '

$ws.Range("D16").Value = '64.925.01'
$ws.Range("E16").Value = '  +1.58%  '

$ws.Range("E17").Value = '  +2.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '525.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +10.50%  '

$ws.Range("D19").Value = '3.138.13'
$ws.Range("E19").Value = '  +3.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.700'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.79%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.74%  '

$ws.Range("E26").Value = '  -0.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.65%  '

$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.08%  '

$ws.Range("E33").Value = '  +4.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '562.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +14.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.75%  '

$ws.Range("E37").Value = '  +9.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0813'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.58%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.65%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.058.80'
$ws.Range("E41").Value = '  +6.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.121'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.256'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.31%  '

$ws.Range("E45").Value = '  +8.65%  '

$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.93%  '

$ws.Range("D49").Value = '0.0₃0523'
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("E50").Value = '  +2.40%  '

$ws.Range("E51").Value = '  +4.16%  '
